# refactor(xlsx-loader): change variables structure, add 'event' column
#
# Both the "pages" and "events" sheets get a new "event" column inserted
# right after "layer" (i.e. becomes column C, pushing key/label/value/
# description/example one column to the right).
#
# - "pages": every data row keeps its "key/label/value" triple; the new
#   "event" column is filled with "page-load" (rows whose layer is "page")
#   or "page-other" (rows whose layer is "other").
# - "events": the old data modelled "event" as just another key/value pair
#   (key="event", value=<event name>) on its own row. That row collapses
#   into the new "event" column on the record's header row, so each record
#   shrinks from 3 sub-rows to 2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "pages"
# ---------------------------------------------------------------
$pages = $wb.Worksheets.Item("pages")

$pages.Columns("C").Insert()

$pages.Range("C1").Value = "event"
$pages.Range("C2").Value = "page-load"
$pages.Range("C5").Value = "page-other"
$pages.Range("C6").Value = "page-load"
$pages.Range("C7").Value = "page-other"
$pages.Range("C8").Value = "page-load"

$pages.Activate() | Out-Null
$pages.Range("C8").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "events"
# ---------------------------------------------------------------
$events = $wb.Worksheets.Item("events")

$events.Columns("C").Insert()

$events.Range("C1").Value = "event"
$events.Range("C2").Value = "search"
$events.Range("C6").Value = "formSubmit"
$events.Range("C9").Value = "formError"

# The old key="event" rows are now redundant - remove them bottom-up so
# earlier row numbers stay valid while deleting.
$events.Rows("10").Delete()
$events.Rows("7").Delete()
$events.Rows("3").Delete()

$events.Activate() | Out-Null
$events.Range("C10").Select() | Out-Null
